# modify SLG building config
# The Prefab column (D) now points at the new "Prefabs/Object/..." assets
# instead of the old "COC_Resources/animation/building/.../1.prefab" paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 4).Value = "Prefabs/Object/Town_1_1"
$ws.Cells.Item(9, 4).Value  = "Prefabs/Object/Tower_1_1"
$ws.Cells.Item(6, 4).Value  = "Prefabs/Object/Item_hourse_1_1"
$ws.Cells.Item(8, 4).Value  = "Prefabs/Object/MagicHourse_1_1"
$ws.Cells.Item(7, 4).Value  = "Prefabs/Object/League_1_1"
$ws.Cells.Item(5, 4).Value  = "Prefabs/Object/GoldMine_1_1"
$ws.Cells.Item(4, 4).Value  = "Prefabs/Object/Camp_1_1"
$ws.Cells.Item(3, 4).Value  = "Prefabs/Object/Arena_1_1"
$ws.Cells.Item(2, 4).Value  = "Prefabs/Object/Altar_1_1"

# Move the selection cursor to match the saved workbook state.
$ws.Range("D11").Select()
